# Adding "coodUncM" (coordinate uncertainty in meters) column, inserted
# between vUTM (U) and numPlants (previously V, now W).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the existing V1 comment (numPlants header note) before the column
# shift so we can re-home it on the cell it describes (now W1).
$oldComment = $ws.Range("V1").Comment
$oldCommentText = $oldComment.Text()

# Insert a new column at V; everything from V onward (numPlants,
# occRemarks, collector, idBy, assOcc, assTaxa, dataEntryRemarks) shifts
# one column to the right.
$ws.Columns("V").Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 22).Value = "coodUncM"

# Move the "numPlants" comment from its old location (V1) to its new one
# (W1), preserving its text.
$oldComment.Delete()
$ws.Range("W1").AddComment($oldCommentText) | Out-Null

# Update the active selection to match the saved workbook state.
$ws.Range("V3").Select()
